# feat: parse `Edm.Time` when excel data is text
#
# Cell H3 previously held a numeric Excel time value (0.6666... with a
# "h:mm:ss" style). It is changed to hold the literal text "16:00:00"
# instead (so that downstream code can parse the raw Edm.Time string).
# Setting NumberFormat to "@" (Text) first forces Excel to store the
# value as text rather than re-interpreting "16:00:00" as a time serial.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("H3")
$cell.NumberFormat = "@"
$cell.Value = "16:00:00"

# Move the active selection, matching the author's saved cursor position.
[void]$ws.Range("H13").Select()
